$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Images")

# Header E2: "erode_c4_unroll" -> "erode_c4_pragma" (unused kernel removed)
$ws.Range("E2").Value = "erode_c4_pragma"

# Widen column E slightly to fit the new header text
$ws.Range("E1").EntireColumn.ColumnWidth = 15.25

# New note row explaining why *_c4_unroll was dropped
$ws.Cells.Item(41, 2).Value = "_c4_unroll jest w zasadzie identyczne co _c4"

# Move the active selection to reflect where the author was last working
$ws.Range("F22").Select() | Out-Null
